# Auto-generated: apply scheduled-runner price/profit refresh to the
# per-job Leve Profit tables (one per worksheet). Each worksheet has the
# same header layout: A:G leve/item info, H:N = currentAveragePrice,
# currentAveragePriceNQ, currentAveragePriceHQ, LevePriceNQ, LevePriceHQ,
# LeveProfitNQ, LeveProfitHQ.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = [ordered]@{
    "H41" = 1484.8
    "I41" = 2096.5
    "J41" = 567.25
    "K41" = 2096.5
    "L41" = 567.25
    "M41" = -1656.5
    "N41" = -1447.25
    "H53" = 1266.1818
    "I53" = 166.66667
    "K53" = 166.66667
    "M53" = 470.33333
    "H86" = 30916.334
    "J86" = 50000
    "L86" = 50000
    "N86" = -52246
    "H89" = 30916.334
    "J89" = 50000
    "L89" = 250000
    "N89" = -261232
    "H92" = 652.6842
    "I92" = 577.8333
    "J92" = 2000
    "K92" = 577.8333
    "L92" = 2000
    "M92" = 670.1667
    "N92" = -4496
    "H98" = 1441.0526
    "I98" = 1024
    "J98" = 3665.3333
    "K98" = 1024
    "L98" = 3665.3333
    "M98" = 474
    "N98" = -6661.3333
    "H122" = 1441.0526
    "I122" = 1024
    "J122" = 3665.3333
    "K122" = 3072
    "L122" = 10995.9999
    "M122" = -622
    "N122" = -15895.9999
    "H128" = 90000
    "J128" = 90000
    "L128" = 90000
    "N128" = -99960
    "H132" = 2793.8667
    "I132" = 2793.8667
    "K132" = 8381.6001
    "M132" = -5851.6001
    "H139" = 104748.75
    "J139" = 104748.75
    "L139" = 104748.75
    "N139" = -115028.75
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("ARM")
$updates = [ordered]@{
    "H25" = 6000
    "I25" = 6000
    "J25" = 0
    "K25" = 6000
    "L25" = 0
    "M25" = -5598
    "H45" = 6723.5
    "J45" = 4934.75
    "L45" = 4934.75
    "N45" = -5688.75
    "H132" = 2327310.8
    "I132" = 2858663.2
    "J132" = 2643.625
    "K132" = 8575989.600000001
    "L132" = 7930.875
    "M132" = -8573459.600000001
    "N132" = -12990.875
    "H141" = 81986.25
    "J141" = 81986.25
    "L141" = 81986.25
    "N141" = -92346.25
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
foreach ($ref in @("N25")) {
    $ws.Range($ref).ClearContents()
}

$ws = $wb.Worksheets.Item("BSM")
$updates = [ordered]@{
    "H15" = 0
    "J15" = 0
    "L15" = 0
    "H86" = 3311.7273
    "J86" = 3709.111
    "L86" = 3709.111
    "N86" = -5955.111
    "H89" = 3311.7273
    "J89" = 3709.111
    "L89" = 18545.555
    "N89" = -29777.555
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
foreach ($ref in @("N15")) {
    $ws.Range($ref).ClearContents()
}

$ws = $wb.Worksheets.Item("CRP")
$updates = [ordered]@{
    "H6" = 3000
    "J6" = 4999
    "L6" = 4999
    "N6" = -5225
    "H16" = 1365.1333
    "J16" = 1132.6666
    "L16" = 1132.6666
    "N16" = -1706.6666
    "H19" = 2037.125
    "I19" = 1849.5
    "K19" = 1849.5
    "M19" = -1679.5
    "H24" = 2037.125
    "I24" = 1849.5
    "K24" = 1849.5
    "M24" = -1679.5
    "H113" = 1365.1333
    "J113" = 1132.6666
    "L113" = 1132.6666
    "N113" = -5472.6666
    "H122" = 1605.4706
    "J122" = 795
    "L122" = 2385
    "N122" = -7285
    "H132" = 58826530
    "I132" = 76925770
    "K132" = 230777310
    "M132" = -230774780
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("CUL")
$updates = [ordered]@{
    "H39" = 700.25
    "I39" = 700.25
    "J39" = 0
    "K39" = 2100.75
    "L39" = 0
    "M39" = -1806.75
    "H122" = 1420
    "J122" = 1800
    "L122" = 16200
    "N122" = -21100
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
foreach ($ref in @("N39")) {
    $ws.Range($ref).ClearContents()
}

$ws = $wb.Worksheets.Item("GSM")
$updates = [ordered]@{
    "H80" = 3079.5386
    "J80" = 2769.75
    "L80" = 2769.75
    "N80" = -4765.75
    "H83" = 3079.5386
    "J83" = 2769.75
    "L83" = 13848.75
    "N83" = -23832.75
    "H102" = 2430.524
    "I102" = 2165.842
    "K102" = 2165.842
    "M102" = -543.8420000000001
    "H113" = 63068.53
    "I113" = 75940.36
    "K113" = 75940.36
    "M113" = -73770.36
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("LTW")
$updates = [ordered]@{
    "H22" = 3095.0715
    "I22" = 3155.5557
    "K22" = 3155.5557
    "M22" = -2860.5557
    "H27" = 3095.0715
    "I27" = 3155.5557
    "K27" = 3155.5557
    "M27" = -3048.5557
    "H40" = 4147.7334
    "I40" = 4147.7334
    "K40" = 4147.7334
    "M40" = -4011.7334
    "H46" = 2564.7273
    "I46" = 2512.4443
    "K46" = 2512.4443
    "M46" = -2324.4443
    "H55" = 573.44446
    "I55" = 432
    "J55" = 750.25
    "K55" = 432
    "L55" = 750.25
    "M55" = -259
    "N55" = -1096.25
    "H61" = 4707.1665
    "I61" = 4952.8184
    "K61" = 4952.8184
    "M61" = -4750.8184
    "H113" = 4707.1665
    "I113" = 4952.8184
    "K113" = 4952.8184
    "M113" = -2782.8184
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("WVR")
$updates = [ordered]@{
    "H41" = 32496.334
    "J41" = 33195.8
    "L41" = 33195.8
    "N41" = -33975.8
    "H74" = 77812.5
    "J74" = 80626
    "L74" = 80626
    "N74" = -82498
    "H77" = 77812.5
    "J77" = 80626
    "L77" = 241878
    "N77" = -251238
    "H122" = 1294.08
    "I122" = 1170.2222
    "J122" = 1612.5714
    "K122" = 3510.6666
    "L122" = 4837.7142
    "M122" = -1060.6666
    "N122" = -9737.7142
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
